$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing rows 2-10 (columns A-F) ----
$existing = @{
    2  = @{ A = 0; B = "NSE:AUTOAXLES";  C = "NSE:ARTEMISMED"; D = "";            E = "NSE:LT";   F = "NSE:CHAMBLFERT" }
    3  = @{ A = 1; B = "NSE:BLUEDART";   C = "NSE:BBL";        D = "";            E = "NSE:NYKAA"; F = "" }
    4  = @{ A = 2; B = "NSE:BOROLTD";    C = "NSE:CAPTRUST";   D = "";            E = "";           F = "" }
    5  = @{ A = 3; B = "NSE:CHAMBLFERT"; C = "NSE:CHOLAHLDNG"; D = "";            E = "";           F = "" }
    6  = @{ A = 4; B = "NSE:CREATIVEYE"; C = "NSE:EASEMYTRIP"; D = "";            E = "";           F = "" }
    7  = @{ A = 5; B = "NSE:DEEPINDS";   C = "NSE:FINEORG";    D = "";            E = "";           F = "" }
    8  = @{ A = 6; B = "NSE:EXCELINDUS"; C = "NSE:GODREJAGRO"; D = "";            E = "";           F = "" }
    9  = @{ A = 7; B = "NSE:FMGOETZE";   C = "NSE:GTL";        D = "";            E = "";           F = "" }
    10 = @{ A = 8; B = "NSE:KRBL";       C = "NSE:HINDUNILVR"; D = "";            E = "";           F = "" }
}

foreach ($r in 2..10) {
    $row = $existing[$r]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}

# ---- Add new rows 11-17 (columns A-F) ----
$newRows = @{
    11 = @{ A = 9;  B = "NSE:ORIENTALTL"; C = "NSE:HTMEDIA" }
    12 = @{ A = 10; B = "";               C = "NSE:INGERRAND" }
    13 = @{ A = 11; B = "";               C = "NSE:KAYNES" }
    14 = @{ A = 12; B = "";               C = "NSE:NESTLEIND" }
    15 = @{ A = 13; B = "";               C = "NSE:NLCINDIA" }
    16 = @{ A = 14; B = "";               C = "NSE:ONMOBILE" }
    17 = @{ A = 15; B = "";               C = "NSE:PRECWIRE" }
}

foreach ($r in 11..17) {
    $row = $newRows[$r]

    # Copy column-A formatting from row 10 so the new index cells keep
    # the same style (s="1": bold, centered, bordered) as the rest of
    # the A column.
    $ws.Cells.Item(10, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
}

$excel.CutCopyMode = $false
